$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (price/volume refresh + two re-ranked rows).
# Values are set with a leading apostrophe to force text, matching the sheet's
# existing inline-string cells (prices like "1.000" / "30.335.34" would otherwise
# be auto-coerced to numbers by Excel and lose formatting).
$ws.Range('D2').Value = "'30.335.34"
$ws.Range('E2').Value = "'  +0.51%  "
$ws.Range('D3').Value = "'1.869.52"
$ws.Range('E3').Value = "'  +0.03%  "
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = "'  -0.12%  "
$ws.Range('D5').Value = "'235.79"
$ws.Range('E5').Value = "'  +0.84%  "
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = "'  -0.10%  "
$ws.Range('E7').Value = "'  -0.03%  "
$ws.Range('D8').Value = "'0.2849"
$ws.Range('E8').Value = "'  +0.69%  "
$ws.Range('D9').Value = "'0.06549"
$ws.Range('E9').Value = "'  -0.58%  "
$ws.Range('D10').Value = "'21.54"
$ws.Range('E10').Value = "'  +6.48%  "
$ws.Range('D11').Value = "'0.07883"
$ws.Range('E11').Value = "'  +1.32%  "
$ws.Range('D12').Value = "'98.13"
$ws.Range('E12').Value = "'  +1.32%  "
$ws.Range('D13').Value = "'1.871.57"
$ws.Range('E13').Value = "'  -0.84%  "
$ws.Range('D14').Value = "'5.105"
$ws.Range('E14').Value = "'  +0.66%  "
$ws.Range('D15').Value = "'0.6765"
$ws.Range('E15').Value = "'  +1.06%  "
$ws.Range('D16').Value = "'277.55"
$ws.Range('E16').Value = "'  -2.95%  "
$ws.Range('D17').Value = "'30.328.99"
$ws.Range('E17').Value = "'  +0.37%  "
$ws.Range('E18').Value = "'  -0.04%  "
$ws.Range('E19').Value = "'  +1.33%  "
$ws.Range('D20').Value = "'5.469"
$ws.Range('E20').Value = "'  +1.84%  "
$ws.Range('D21').Value = "'2.117.20"
$ws.Range('E21').Value = "'  -0.65%  "
$ws.Range('D22').Value = "'0.000007313"
$ws.Range('E22').Value = "'  +0.94%  "
$ws.Range('E23').Value = "'  -0.17%  "
$ws.Range('D24').Value = "'6.154"
$ws.Range('E24').Value = "'  +0.00%  "
$ws.Range('D25').Value = "'165.49"
$ws.Range('E25').Value = "'  -1.13%  "
$ws.Range('D26').Value = "'9.149"
$ws.Range('E26').Value = "'  -2.00%  "
$ws.Range('D27').Value = "'19.13"
$ws.Range('E27').Value = "'  +0.20%  "
$ws.Range('D28').Value = "'1.935"
$ws.Range('E28').Value = "'  -1.30%  "
$ws.Range('E29').Value = "'  +0.64%  "
$ws.Range('D30').Value = "'0.09640"
$ws.Range('E30').Value = "'  -0.06%  "
$ws.Range('D31').Value = "'4.393"
$ws.Range('E31').Value = "'  +0.28%  "
$ws.Range('E32').Value = "'  +0.71%  "
$ws.Range('D33').Value = "'4.098"
$ws.Range('E33').Value = "'  -0.21%  "
$ws.Range('D34').Value = "'0.04701"
$ws.Range('E34').Value = "'  +0.73%  "
$ws.Range('E35').Value = "'  +3.62%  "
$ws.Range('D36').Value = "'0.7064"
$ws.Range('E36').Value = "'  +0.61%  "
$ws.Range('D37').Value = "'2.728"
$ws.Range('E37').Value = "'  +0.30%  "
$ws.Range('D38').Value = "'0.01860"
$ws.Range('E38').Value = "'  -0.25%  "
$ws.Range('D39').Value = "'6.337"
$ws.Range('E39').Value = "'  -1.38%  "
$ws.Range('D40').Value = "'2.539"
$ws.Range('E40').Value = "'  +0.70%  "
$ws.Range('D41').Value = "'74.18"
$ws.Range('E41').Value = "'  +3.61%  "
$ws.Range('D42').Value = "'1.956"
$ws.Range('E42').Value = "'  +0.63%  "
$ws.Range('D43').Value = "'0.8501"
$ws.Range('E43').Value = "'  -1.03%  "
$ws.Range('D44').Value = "'0.4190"
$ws.Range('E44').Value = "'  +0.27%  "
$ws.Range('B45').Value = "'Quant"
$ws.Range('C45').Value = "'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range('D45').Value = "'103.98"
$ws.Range('E45').Value = "'  +0.96%  "
$ws.Range('B46').Value = "'PaxDollar"
$ws.Range('C46').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('D46').Value = "'1.000"
$ws.Range('E46').Value = "'  -0.07%  "
$ws.Range('D47').Value = "'7.205"
$ws.Range('E47').Value = "'  +0.37%  "
$ws.Range('D48').Value = "'9.226"
$ws.Range('E48').Value = "'  +0.61%  "
$ws.Range('D49').Value = "'938.03"
$ws.Range('E49').Value = "'  -5.11%  "
$ws.Range('D50').Value = "'34.20"
$ws.Range('E50').Value = "'  +0.93%  "
$ws.Range('B51').Value = "'Algorand"
$ws.Range('C51').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('D51').Value = "'0.1123"
$ws.Range('E51').Value = "'  -2.05%  "
